{"js": "const replacements = [\n  [\"2024-07-19 Friday\", \"2024-07-20 Saturday\"],\n  [\"38+33=71\", \"95-56=39\"],\n  [\"41+13=54\", \"75-50=25\"],\n  [\"98-21=77\", \"18+61=79\"],\n  [\"77-8=69\", \"31+10=41\"],\n  [\"29+70=99\", \"6+27=33\"],\n  [\"37+54=91\", \"92-42=50\"],\n  [\"50+11=61\", \"72-0=72\"],\n  [\"1+17=18\", \"58+38=96\"],\n  [\"53-27=26\", \"98-1=97\"],\n  [\"6+92=98\", \"19+71=90\"],\n  [\"26+43=69\", \"37-19=18\"],\n  [\"16+50=66\", \"56+7=63\"],\n  [\"60-38=22\", \"89-32=57\"],\n  [\"21+77=98\", \"16+26=42\"],\n  [\"3+84=87\", \"43+17=60\"],\n  [\"89-60=29\", \"65-49=16\"],\n  [\"87-37=50\", \"72+3=75\"],\n  [\"74-72=2\", \"47+4=51\"],\n  [\"93-51=42\", \"59-38=21\"],\n  [\"39-10=29\", \"24+70=94\"],\n  [\"22+29=51\", \"63-10=53\"],\n  [\"42-19=23\", \"36+55=91\"],\n  [\"79-8=71\", \"15+82=97\"],\n  [\"82+1=83\", \"91-7=84\"],\n  [\"28+38=66\", \"38+10=48\"],\n  [\"11+80=91\", \"7+26=33\"],\n  [\"94-21=73\", \"28+8=36\"],\n  [\"30+54=84\", \"74-4=70\"],\n  [\"54-2=52\", \"73+6=79\"],\n  [\"75-35=40\", \"88-41=47\"],\n  [\"35+35=70\", \"51+19=70\"],\n  [\"50-23=27\", \"47-45=2\"],\n  [\"8+80=88\", \"87+8=95\"],\n  [\"21-14=7\", \"82-77=5\"],\n  [\"34-6=28\", \"81-14=67\"],\n  [\"4+34=38\", \"65+30=95\"],\n  [\"7+64=71\", \"38-1=37\"],\n  [\"72+8=80\", \"93-25=68\"],\n  [\"54+42=96\", \"66-22=44\"],\n  [\"68+17=85\", \"99-15=84\"],\n  [\"91-32=59\", \"48-31=17\"],\n  [\"55+41=96\", \"95-38=57\"],\n  [\"29+2=31\", \"44-15=29\"],\n  [\"47-10=37\", \"94-93=1\"],\n  [\"66-60=6\", \"15+33=48\"],\n  [\"82-12=70\", \"46+3=49\"],\n  [\"96-93=3\", \"26+37=63\"],\n  [\"67-23=44\", \"59+1=60\"],\n  [\"31+11=42\", \"74-2=72\"],\n  [\"94-42=52\", \"74+9=83\"],\n  [\"46+32=78\", \"16+59=75\"],\n  [\"21+68=89\", \"29-24=5\"],\n  [\"59-8=51\", \"18-5=13\"],\n  [\"92-38=54\", \"41+45=86\"],\n  [\"53+28=81\", \"9+59=68\"],\n  [\"67-48=19\", \"98-70=28\"],\n  [\"10+22=32\", \"19+50=69\"],\n  [\"48+4=52\", \"8+66=74\"],\n  [\"74+16=90\", \"51-18=33\"],\n  [\"83-54=29\", \"85-41=44\"],\n  [\"13+19=32\", \"88-88=0\"],\n  [\"87-82=5\", \"0+40=40\"],\n  [\"65-12=53\", \"16+8=24\"],\n  [\"14+35=49\", \"18+33=51\"],\n  [\"55+5=60\", \"78-20=58\"],\n  [\"79-46=33\", \"0+0=0\"],\n  [\"71-0=71\", \"77-48=29\"],\n  [\"80-62=18\", \"2+89=91\"],\n  [\"35-25=10\", \"72+1=73\"],\n  [\"89-13=76\", \"82-14=68\"],\n  [\"53+30=83\", \"2+41=43\"],\n  [\"51-28=23\", \"63+8=71\"],\n  [\"40+30=70\", \"49-5=44\"],\n  [\"70-57=13\", \"58-21=37\"],\n  [\"13+39=52\", \"22+32=54\"],\n  [\"27+54=81\", \"21+64=85\"],\n  [\"26-21=5\", \"73+0=73\"],\n  [\"23+71=94\", \"8+51=59\"],\n  [\"16+1=17\", \"84-18=66\"],\n  [\"42-41=1\", \"41-14=27\"],\n  [\"36-7=29\", \"95-67=28\"],\n  [\"27+30=57\", \"51+30=81\"],\n  [\"93-30=63\", \"17+69=86\"],\n  [\"88-26=62\", \"15+69=84\"],\n  [\"93-29=64\", \"25+58=83\"],\n  [\"86-11=75\", \"51+26=77\"],\n  [\"82-56=26\", \"19+73=92\"],\n  [\"22+6=28\", \"48+1=49\"],\n  [\"94+3=97\", \"51+43=94\"],\n  [\"95-11=84\", \"63-14=49\"],\n  [\"20+5=25\", \"35-20=15\"],\n  [\"18+22=40\", \"88-5=83\"],\n  [\"76-63=13\", \"39+51=90\"],\n  [\"47-38=9\", \"51-46=5\"],\n  [\"59+27=86\", \"0+71=71\"],\n  [\"15+45=60\", \"68-45=23\"],\n  [\"33+48=81\", \"70-53=17\"],\n  [\"20+44=64\", \"54-47=7\"],\n  [\"99-55=44\", \"13+69=82\"],\n  [\"35-23=12\", \"40+59=99\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"2024-07-19 Friday\", \"2024-07-20 Saturday\")\n  ,@(\"38+33=71\", \"95-56=39\")\n  ,@(\"41+13=54\", \"75-50=25\")\n  ,@(\"98-21=77\", \"18+61=79\")\n  ,@(\"77-8=69\", \"31+10=41\")\n  ,@(\"29+70=99\", \"6+27=33\")\n  ,@(\"37+54=91\", \"92-42=50\")\n  ,@(\"50+11=61\", \"72-0=72\")\n  ,@(\"1+17=18\", \"58+38=96\")\n  ,@(\"53-27=26\", \"98-1=97\")\n  ,@(\"6+92=98\", \"19+71=90\")\n  ,@(\"26+43=69\", \"37-19=18\")\n  ,@(\"16+50=66\", \"56+7=63\")\n  ,@(\"60-38=22\", \"89-32=57\")\n  ,@(\"21+77=98\", \"16+26=42\")\n  ,@(\"3+84=87\", \"43+17=60\")\n  ,@(\"89-60=29\", \"65-49=16\")\n  ,@(\"87-37=50\", \"72+3=75\")\n  ,@(\"74-72=2\", \"47+4=51\")\n  ,@(\"93-51=42\", \"59-38=21\")\n  ,@(\"39-10=29\", \"24+70=94\")\n  ,@(\"22+29=51\", \"63-10=53\")\n  ,@(\"42-19=23\", \"36+55=91\")\n  ,@(\"79-8=71\", \"15+82=97\")\n  ,@(\"82+1=83\", \"91-7=84\")\n  ,@(\"28+38=66\", \"38+10=48\")\n  ,@(\"11+80=91\", \"7+26=33\")\n  ,@(\"94-21=73\", \"28+8=36\")\n  ,@(\"30+54=84\", \"74-4=70\")\n  ,@(\"54-2=52\", \"73+6=79\")\n  ,@(\"75-35=40\", \"88-41=47\")\n  ,@(\"35+35=70\", \"51+19=70\")\n  ,@(\"50-23=27\", \"47-45=2\")\n  ,@(\"8+80=88\", \"87+8=95\")\n  ,@(\"21-14=7\", \"82-77=5\")\n  ,@(\"34-6=28\", \"81-14=67\")\n  ,@(\"4+34=38\", \"65+30=95\")\n  ,@(\"7+64=71\", \"38-1=37\")\n  ,@(\"72+8=80\", \"93-25=68\")\n  ,@(\"54+42=96\", \"66-22=44\")\n  ,@(\"68+17=85\", \"99-15=84\")\n  ,@(\"91-32=59\", \"48-31=17\")\n  ,@(\"55+41=96\", \"95-38=57\")\n  ,@(\"29+2=31\", \"44-15=29\")\n  ,@(\"47-10=37\", \"94-93=1\")\n  ,@(\"66-60=6\", \"15+33=48\")\n  ,@(\"82-12=70\", \"46+3=49\")\n  ,@(\"96-93=3\", \"26+37=63\")\n  ,@(\"67-23=44\", \"59+1=60\")\n  ,@(\"31+11=42\", \"74-2=72\")\n  ,@(\"94-42=52\", \"74+9=83\")\n  ,@(\"46+32=78\", \"16+59=75\")\n  ,@(\"21+68=89\", \"29-24=5\")\n  ,@(\"59-8=51\", \"18-5=13\")\n  ,@(\"92-38=54\", \"41+45=86\")\n  ,@(\"53+28=81\", \"9+59=68\")\n  ,@(\"67-48=19\", \"98-70=28\")\n  ,@(\"10+22=32\", \"19+50=69\")\n  ,@(\"48+4=52\", \"8+66=74\")\n  ,@(\"74+16=90\", \"51-18=33\")\n  ,@(\"83-54=29\", \"85-41=44\")\n  ,@(\"13+19=32\", \"88-88=0\")\n  ,@(\"87-82=5\", \"0+40=40\")\n  ,@(\"65-12=53\", \"16+8=24\")\n  ,@(\"14+35=49\", \"18+33=51\")\n  ,@(\"55+5=60\", \"78-20=58\")\n  ,@(\"79-46=33\", \"0+0=0\")\n  ,@(\"71-0=71\", \"77-48=29\")\n  ,@(\"80-62=18\", \"2+89=91\")\n  ,@(\"35-25=10\", \"72+1=73\")\n  ,@(\"89-13=76\", \"82-14=68\")\n  ,@(\"53+30=83\", \"2+41=43\")\n  ,@(\"51-28=23\", \"63+8=71\")\n  ,@(\"40+30=70\", \"49-5=44\")\n  ,@(\"70-57=13\", \"58-21=37\")\n  ,@(\"13+39=52\", \"22+32=54\")\n  ,@(\"27+54=81\", \"21+64=85\")\n  ,@(\"26-21=5\", \"73+0=73\")\n  ,@(\"23+71=94\", \"8+51=59\")\n  ,@(\"16+1=17\", \"84-18=66\")\n  ,@(\"42-41=1\", \"41-14=27\")\n  ,@(\"36-7=29\", \"95-67=28\")\n  ,@(\"27+30=57\", \"51+30=81\")\n  ,@(\"93-30=63\", \"17+69=86\")\n  ,@(\"88-26=62\", \"15+69=84\")\n  ,@(\"93-29=64\", \"25+58=83\")\n  ,@(\"86-11=75\", \"51+26=77\")\n  ,@(\"82-56=26\", \"19+73=92\")\n  ,@(\"22+6=28\", \"48+1=49\")\n  ,@(\"94+3=97\", \"51+43=94\")\n  ,@(\"95-11=84\", \"63-14=49\")\n  ,@(\"20+5=25\", \"35-20=15\")\n  ,@(\"18+22=40\", \"88-5=83\")\n  ,@(\"76-63=13\", \"39+51=90\")\n  ,@(\"47-38=9\", \"51-46=5\")\n  ,@(\"59+27=86\", \"0+71=71\")\n  ,@(\"15+45=60\", \"68-45=23\")\n  ,@(\"33+48=81\", \"70-53=17\")\n  ,@(\"20+44=64\", \"54-47=7\")\n  ,@(\"99-55=44\", \"13+69=82\")\n  ,@(\"35-23=12\", \"40+59=99\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $r = $d.Content\n  $found = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}"}
